$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current (pre-edit) values for rows 4, 5, 6 in columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")
$old = @{}
foreach ($row in 4..6) {
    $old[$row] = @{}
    foreach ($col in $cols) {
        $old[$row][$col] = $ws.Range("$col$row").Value2
    }
}

# Row 4 takes old Row 5 values; Row 5 takes old Row 6 values; Row 6 takes old Row 4 values
$mapping = @{ 4 = 5; 5 = 6; 6 = 4 }

foreach ($row in 4..6) {
    $srcRow = $mapping[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value2 = $old[$srcRow][$col]
    }
}
